# Update Excel file with latest predictions
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Home win
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Home win")
$rows = @(
  @("20-12-2024 19:45", "ENGLAND", "LEAGUE ONE", "Stockport County - Peterborough", 70, 1.76),
  @("20-12-2024 13:00", "ISRAEL", "LIGA LEUMIT", "Hapoel Kfar Shalem - Hapoel Ramat HaSharon", 80, 1.91),
  @("20-12-2024 13:00", "WORLD", "AFF CHAMPIONSHIP", "Malaysia - Singapore", 80, 1.7),
  @("21-12-2024 15:00", "ENGLAND", "LEAGUE TWO", "Notts County - Bradford", 73.3, 2.2),
  @("21-12-2024 15:00", "SCOTLAND", "CHAMPIONSHIP", "Livingston - Ayr Utd", 80, 2),
  @("21-12-2024 15:00", "SCOTLAND", "LEAGUE ONE", "Kelty Hearts - Queen Of The South", 71.7, 2.3),
  @("21-12-2024 19:00", "BELGIUM", "CHALLENGER PRO LEAGUE", "RAAL La Louvière - RWDM", 70, 2.2),
  @("21-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE - NORTH", "South Shields - Marine", 73.3, 2.15),
  @("21-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE - SOUTH", "Bath City - Aveley", 73.3, 1.95),
  @("21-12-2024 13:00", "SPAIN", "SEGUNDA DIVISIÓN", "Eibar - Granada CF", 80, 2.37)
)
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $ws.Cells.Item($r,6).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet: Draw
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Draw")
$rows = @(
  @("21-12-2024 12:30", "ENGLAND", "LEAGUE ONE", "Lincoln - Reading", 66.7, 3.7),
  @("21-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE - SOUTH", "Weston-super-Mare - Eastbourne Borough", 60, 3.25),
  @("21-12-2024 15:30", "GERMANY", "3. LIGA", "Alemannia Aachen - SV Wehen", 60, 3.3),
  @("21-12-2024 16:30", "ITALY", "SERIE C - GIRONE C", "Casertana - Latina", 66.7, 3),
  @("21-12-2024 11:00", "PORTUGAL", "SEGUNDA LIGA", "Vizela - FC Porto B", 60, 3.35),
  @("21-12-2024 17:30", "SPAIN", "SEGUNDA DIVISIÓN", "Zaragoza - Racing Ferrol", 60, 3.5)
)
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $ws.Cells.Item($r,6).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet: Btts
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Btts")
$rows = @(
  @("20-12-2024 19:45", "ENGLAND", "LEAGUE ONE", "Cambridge United - Huddersfield", 76, 1.8),
  @("20-12-2024 19:30", "ITALY", "SERIE C - GIRONE B", "Pontedera - Legnago Salus", 88, 1.85),
  @("21-12-2024 12:30", "ENGLAND", "CHAMPIONSHIP", "Hull City - Swansea", 76.7, 1.73),
  @("21-12-2024 12:30", "ENGLAND", "LEAGUE ONE", "Lincoln - Reading", 80, 1.75),
  @("21-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE", "Braintree - Oldham", 76, 1.8),
  @("21-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE", "Eastleigh - Maidenhead", 80, 1.75),
  @("21-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE", "Solihull Moors - AFC Fylde", 80, 1.75),
  @("21-12-2024 20:30", "PORTUGAL", "PRIMEIRA LIGA", "Moreirense - FC Porto", 76.7, 2.1),
  @("21-12-2024 15:00", "SCOTLAND", "LEAGUE TWO", "Peterhead - Bonnyrigg Rose Athletic", 78.3, 1.73),
  @("21-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE - NORTH", "Scarborough Athletic - Peterborough Sports", 83.3, 1.85),
  @("21-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE - SOUTH", "Chelmsford City - Weymouth", 78.3, 1.73),
  @("21-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE - SOUTH", "Worthing - Enfield Town", 78.3, 1.73),
  @("21-12-2024 15:30", "GERMANY", "3. LIGA", "Alemannia Aachen - SV Wehen", 88, 1.7),
  @("21-12-2024 14:00", "ITALY", "SERIE C - GIRONE A", "Clodiense - Lecco", 76, 1.9),
  @("21-12-2024 14:00", "ITALY", "SERIE C - GIRONE A", "Lumezzane - Virtus Verona", 80, 1.95),
  @("21-12-2024 15:00", "MOROCCO", "BOTOLA PRO", "Riadi Salmi - Moghreb Tetouan", 76.7, 2.2),
  @("21-12-2024 18:30", "SPAIN", "PRIMERA DIVISIÓN RFEF - GROUP 2", "Mérida AD - Fuenlabrada", 76.7, 2),
  @("21-12-2024 13:00", "TURKEY", "SÜPER LIG", "Sivasspor - Samsunspor", 80, 1.73)
)
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $ws.Cells.Item($r,6).Value = $row[5]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Sheet: Over_Under
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Over_Under")
$rows = @(
  @("20-12-2024 19:45", "ITALY", "SERIE A", "Verona - AC Milan", 85, 1.73, 55, 2.75),
  @("20-12-2024 19:00", "NETHERLANDS", "EERSTE DIVISIE", "Dordrecht - FC Eindhoven", 80, 1.48, 80, 2.2),
  @("20-12-2024 19:00", "BELGIUM", "CHALLENGER PRO LEAGUE", "Club Brugge II - Zulte Waregem", 85, 1.7, 65, 2.75),
  @("20-12-2024 13:00", "WORLD", "AFF CHAMPIONSHIP", "Malaysia - Singapore", 75, 1.65, 60, 2.63),
  @("21-12-2024 20:00", "NETHERLANDS", "EREDIVISIE", "Heracles - Groningen", 75, 1.95, 60, 3.4),
  @("21-12-2024 15:00", "SCOTLAND", "LEAGUE TWO", "Stirling Albion - Stranraer", 80, 1.9, 50, 3.25),
  @("21-12-2024 15:00", "ENGLAND", "NATIONAL LEAGUE - NORTH", "Chorley - Rushall Olympic", 95, 1.75, 30, 3),
  @("21-12-2024 16:30", "ITALY", "SERIE C - GIRONE C", "Team Altamura - Crotone", 86.7, 1.73, 46.7, 2.75),
  @("21-12-2024 15:00", "NORTHERN-IRELAND", "PREMIERSHIP", "Loughgall - Ballymena United", 80, 1.8, 70, 2.88)
)
$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r,1).Value = $row[0]
    $ws.Cells.Item($r,2).Value = $row[1]
    $ws.Cells.Item($r,3).Value = $row[2]
    $ws.Cells.Item($r,4).Value = $row[3]
    $ws.Cells.Item($r,5).Value = $row[4]
    $ws.Cells.Item($r,6).Value = $row[5]
    $ws.Cells.Item($r,7).Value = $row[6]
    $ws.Cells.Item($r,8).Value = $row[7]
    $r = $r + 1
}

Write-Host "Update complete"
